# edit.ps1 - apply the "Update document formatting and remove unnecessary
# sections for clarity" change set to the IMSKLK1KT-Sample.docx template.
#
# Summary of the change:
#   * All paragraphs using the custom "Header Style" paragraph style are
#     switched to the built-in "Heading 2" style.
#   * Several section headings get new wording, and a couple of new
#     Heading-2 headings are inserted.
#   * One now-redundant heading paragraph is deleted outright.
#   * The trailing "Publications / review" boilerplate is trimmed out of
#     the big "Data Analysis" paragraph.
#   * The "Heading 2" style definition itself is re-themed (Calibri,
#     navy 000080, 12pt) to match how "Header Style" used to look.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

function Get-ParaByText($text) {
    # Locate the (unique) paragraph whose text matches $text exactly and
    # return the Paragraph object.
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($text, $true, $true, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $text"
    }
    return $rng.Paragraphs(1)
}

function Set-ParaHeading2($text) {
    # Re-style the paragraph whose (current) text equals $text to Heading 2.
    $p = Get-ParaByText($text)
    $p.Style = "Heading 2"
    return $p
}

function Replace-Text($oldText, $newText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, `
                  $false, $newText, 2) | Out-Null
}

function Insert-HeadingAfter($para, $text) {
    # Insert a brand-new Heading-2 paragraph right after $para, with $text
    # as its content. Returns the new Paragraph object.
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $newPara.Range.Text = $text
    $newPara.Style = "Heading 2"
    return $newPara
}

# ---------------------------------------------------------------------
# 1. INTENDED USE -- style only
# ---------------------------------------------------------------------
Set-ParaHeading2("INTENDED USE") | Out-Null

# ---------------------------------------------------------------------
# 2. BACKGROUND -- style + text
# ---------------------------------------------------------------------
Replace-Text "BACKGROUND" "BACKGROUND ON Mouse KLK1 ELISA Kit"
Set-ParaHeading2("BACKGROUND ON Mouse KLK1 ELISA Kit") | Out-Null

# ---------------------------------------------------------------------
# 3. PRINCIPLE OF THE ASSAY -- style only
# ---------------------------------------------------------------------
Set-ParaHeading2("PRINCIPLE OF THE ASSAY") | Out-Null

# ---------------------------------------------------------------------
# 4-6. SPECIFICATION / REAGENTS / MATERIALS REQUIRED BUT NOT PROVIDED
#      -> OVERVIEW / TECHNICAL DETAILS / PREPARATIONS BEFORE ASSAY
#      + two new headings inserted after.
# ---------------------------------------------------------------------
Replace-Text "SPECIFICATION" "OVERVIEW"
Set-ParaHeading2("OVERVIEW") | Out-Null

Replace-Text "REAGENTS" "TECHNICAL DETAILS"
Set-ParaHeading2("TECHNICAL DETAILS") | Out-Null

Replace-Text "MATERIALS REQUIRED BUT NOT PROVIDED" "PREPARATIONS BEFORE ASSAY"
$pPrep = Set-ParaHeading2("PREPARATIONS BEFORE ASSAY")

$pKit = Insert-HeadingAfter $pPrep "KIT COMPONENTS/MATERIALS PROVIDED"
Insert-HeadingAfter $pKit "REQUIRED MATERIALS THAT ARE NOT SUPPLIED" | Out-Null

# ---------------------------------------------------------------------
# 7-8. TYPICAL DATA / TYPICAL STANDARD CURVE
# ---------------------------------------------------------------------
Set-ParaHeading2("TYPICAL DATA") | Out-Null

Replace-Text "TYPICAL STANDARD CURVE" "Mouse KLK1 ELISA Kit STANDARD CURVE EXAMPLE"
Set-ParaHeading2("Mouse KLK1 ELISA Kit STANDARD CURVE EXAMPLE") | Out-Null

# ---------------------------------------------------------------------
# 9. INTRA/INTER ASSAY VARIABILITY -> INTRA/INTER-ASSAY VARIABILITY
# ---------------------------------------------------------------------
Replace-Text "INTRA/INTER ASSAY VARIABILITY" "INTRA/INTER-ASSAY VARIABILITY"
Set-ParaHeading2("INTRA/INTER-ASSAY VARIABILITY") | Out-Null

# ---------------------------------------------------------------------
# 10. REPRODUCIBILITY -- style only
# ---------------------------------------------------------------------
Set-ParaHeading2("REPRODUCIBILITY") | Out-Null

# ---------------------------------------------------------------------
# 11. PROCEDURAL NOTES -> PREPARATION BEFORE THE EXPERIMENT
# ---------------------------------------------------------------------
Replace-Text "PROCEDURAL NOTES" "PREPARATION BEFORE THE EXPERIMENT"
Set-ParaHeading2("PREPARATION BEFORE THE EXPERIMENT") | Out-Null

# ---------------------------------------------------------------------
# 12. Delete the "REAGENT PREPARATION AND STORAGE" heading paragraph.
# ---------------------------------------------------------------------
$pReagentPrep = Get-ParaByText("REAGENT PREPARATION AND STORAGE")
$pReagentPrep.Range.Delete()

# ---------------------------------------------------------------------
# 13. DILUTION OF STANDARD -> DILUTION OF Mouse KLK1 ELISA Kit STANDARD
# ---------------------------------------------------------------------
Replace-Text "DILUTION OF STANDARD" "DILUTION OF Mouse KLK1 ELISA Kit STANDARD"
Set-ParaHeading2("DILUTION OF Mouse KLK1 ELISA Kit STANDARD") | Out-Null

# ---------------------------------------------------------------------
# 14. SAMPLE COLLECTION & STORAGE -> SAMPLE PREPARATION AND STORAGE
#     + two new headings inserted after the following body paragraph.
# ---------------------------------------------------------------------
Replace-Text "SAMPLE COLLECTION & STORAGE" "SAMPLE PREPARATION AND STORAGE"
$pSample = Set-ParaHeading2("SAMPLE PREPARATION AND STORAGE")

$pSampleBody = $pSample.Next()
$pNotes = Insert-HeadingAfter $pSampleBody "SAMPLE COLLECTION NOTES"
Insert-HeadingAfter $pNotes "SAMPLE DILUTION GUIDELINE" | Out-Null

# ---------------------------------------------------------------------
# 15. ASSAY PROCEDURE -- style only
# ---------------------------------------------------------------------
Set-ParaHeading2("ASSAY PROCEDURE") | Out-Null

# ---------------------------------------------------------------------
# 16. DATA ANALYSIS -- style only
# ---------------------------------------------------------------------
Set-ParaHeading2("DATA ANALYSIS") | Out-Null

# ---------------------------------------------------------------------
# 17. Trim the trailing "Publications citing this product" / "Submit a
#     review" boilerplate off the end of the Data Analysis paragraph.
#     (Directly overwrite the sub-range from "PubMed ID" through the end
#     of the paragraph, rather than a Find/Replace of the whole blob, so
#     we do not have to fight the console's mangling of the non-ASCII
#     registered-trademark character that used to sit at the very end.)
# ---------------------------------------------------------------------
$pDataAnalysis = Get-ParaByText("DATA ANALYSIS")
$pDataBody = $pDataAnalysis.Next()

$rngPubMed = $d.Content
$rngPubMed.Find.ClearFormatting()
$okPubMed = $rngPubMed.Find.Execute("PubMed ID", $true, $true, $false, `
                                     $false, $false, $true, 1, $false, "", 0)
if (-not $okPubMed) {
    throw "Text not found: PubMed ID"
}

$newTail = "PubMed ID: 10.1186/s12014-021-09335-9, Proteomics and " + `
    "functional study reveal kallikrein-6 enhances communicating " + `
    "hydrocephalus Visit bosterbio.com/mouse-klk1-picokine-trade-elisa-" + `
    "kit-ek1586-innovative research.html to see all 1 publications. "

$subRange = $d.Range($rngPubMed.Start, $pDataBody.Range.End - 1)
$subRange.Text = $newTail

# ---------------------------------------------------------------------
# 18. DISCLAIMER -- style only
# ---------------------------------------------------------------------
Set-ParaHeading2("DISCLAIMER") | Out-Null

# ---------------------------------------------------------------------
# 19. Re-theme the "Heading 2" style itself: Calibri, navy (000080), 12pt
#     (visually matching how "Header Style" used to render).
# ---------------------------------------------------------------------
$h2 = $d.Styles("Heading 2")
$h2.Font.Name = "Calibri"
$h2.Font.Color = 8388608   # wdColor BGR for RGB 000080 (navy)
$h2.Font.Size = 12

Write-Host "Done."
